# "Update AM slides (#2)" -- three text corrections in the AM Session deck.
$p = $ppt.ActivePresentation

# --- Slide 18: "rule I1 + I1 => I1 +Int I2" -> "rule I1 + I2 => I1 +Int I2"
$s18 = $p.Slides.Item(18)
$sh18 = $s18.Shapes.Item(6)
$run18 = $sh18.TextFrame.TextRange.Paragraphs(4, 1).Runs(10, 1)
$run18.Text = "rule I1 + I2 => I1 +Int I2"

# --- Slide 26: split "constructors (i.e., ... itself)." into two runs and
#     reword "all distinct" -> "configurations that are distinct", coloring
#     the parenthetical with the dk1 theme color.
$s26 = $p.Slides.Item(26)
$sh26 = $s26.Shapes.Item(6)
$para26 = $sh26.TextFrame.TextRange.Paragraphs(1, 1)

$run26a = $para26.Runs(6, 1)
$run26a.Text = "constructors "

$run26b = $para26.Runs(7, 1)
$run26b.Text = "(i.e., configurations that are distinct from each other and only matches itself)."
$run26b.Font.Color.ObjectThemeColor = 1

# --- Slide 27: "(greatest fixpoint - νX.φ)" -> "(greatest fixpoint, νX.φ)"
$s27 = $p.Slides.Item(27)
$sh27 = $s27.Shapes.Item(6)
$run27 = $sh27.TextFrame.TextRange.Paragraphs(7, 1).Runs(1, 1)
$run27.Text = "Either φ holds immediately (greatest fixpoint, νX.φ)"
